# mvp feature completed: AddClients
#
# The "email" column previously held per-row mailto hyperlinks pointing at
# each existing client's address. As part of wiring up the AddClients
# feature, those stale hyperlinked addresses are replaced with the new
# client's plain email address, and the now-unused hyperlinks are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the mailto hyperlinks that used to live on C2:C4 (the cells keep
# their existing "hyperlink-look" style, only the link itself goes away).
$ws.Hyperlinks.Delete()

# New client's email replaces the old per-row addresses in column C.
$newClientEmail = "dumblaymyhit@gmail.com"
$ws.Range("C2").Value = $newClientEmail
$ws.Range("C3").Value = $newClientEmail
$ws.Range("C4").Value = $newClientEmail

# Leave the selection where the user was last working.
$ws.Range("C4").Select()
